# Auto-generated Excel COM-interop script.
# Applies text/label edits across multiple worksheets and renames
# the 'cumcontrol1' sheet to 'cumcontrol', per the target commit diff.

$wb = $excel.ActiveWorkbook

# --- price1 (sheet7) ---
$ws = $wb.Worksheets.Item("price1")
$ws.Range("B2").Value = "I'll hold onto it then, but next time I'm feeling like this you better be ready"
$ws.Range("B3").Value = "alright [lower price] but only because this convo has been different, don't tell anyone"
$ws.Range("B4").Value = "maybe you're not ready for what I just did in this one"
$ws.Range("B5").Value = "I'm only in this mood because of you rn, no guarantee it happens again"
$ws.Range("B6").Value = "that's less than your morning coffee and trust me this hits way harder"
$ws.Range("C2").Value = "SEED. Continue."

# --- price2 (sheet8) ---
$ws = $wb.Worksheets.Item("price2")
$ws.Range("B2").Value = "I'll keep it for now, but I know you'll come back for it"
$ws.Range("B3").Value = "[lower price] because I actually want you to have this one, I'm not making this offer to anyone else"
$ws.Range("B4").Value = "I don't think you're ready for this one honestly, it's a lot"
$ws.Range("B5").Value = "this doesn't happen often and I already chose you to share it with, don't make me regret it"
$ws.Range("B6").Value = "it's literally nothing for what you're about to see, you'll be replaying this one"

# --- discount1 (sheet9) ---
$ws = $wb.Worksheets.Item("discount1")
$ws.Range("B2").Value = "okay I'll keep it then, someone else has been asking and they won't hesitate"
$ws.Range("B3").Value = "fine [lower price] ONLY because I like you, but don't ever ask me this again"
$ws.Range("B4").Value = "the girls who get my best stuff never ask for discounts, they know what it's worth"
$ws.Range("B5").Value = "haha negotiate? this isn't a sale babe, you know exactly what you're getting"

# --- discount2 (sheet10) ---
$ws = $wb.Worksheets.Item("discount2")
$ws.Range("B2").Value = "no problem, I'll save it for someone who really wants it"
$ws.Range("B3").Value = "[lower price] because you've been making me feel a type of way, first and last time"
$ws.Range("B4").Value = "the ones who really want it don't ask for discounts, just saying"
$ws.Range("B5").Value = "a discount? do I look like I'm on sale?"

# --- free1 (sheet11) ---
$ws = $wb.Worksheets.Item("free1")
$ws.Range("B2").Value = "I'll keep it then, but it's got your name on it"
$ws.Range("B3").Value = "I made this because of what YOU said to me, this wasn't random, I did it for you"
$ws.Range("B4").Value = "free? I don't just hand this out, you have to earn the best stuff"
$ws.Range("B5").Value = "I already gave you one for free remember? this one is way crazier, you know I don't do free twice"

# --- free2 (sheet12) ---
$ws = $wb.Worksheets.Item("free2")
$ws.Range("B2").Value = "okay I'll hold onto it for now"
$ws.Range("B3").Value = "I recorded this because of our conversation, specifically because of you, and it took time"
$ws.Range("B4").Value = "the best things aren't free babe, especially not what I just recorded"
$ws.Range("B5").Value = "you already got the free preview, this one is ten times better"

# --- nomoney1 (sheet13) ---
$ws = $wb.Worksheets.Item("nomoney1")
$ws.Range("B2").Value = "it's fine, you're here and that's what counts right now"
$ws.Range("B3").Value = "send whatever you can, even something small, I just need you to see what you made me do"
$ws.Range("B4").Value = "not even [small amount]? I really want you to be the one who sees this"
$ws.Range("B5").Value = "I hear you, that's okay"
$ws.Range("C4").Value = "TEST $3-5. Still no → Step 3."

# --- nomoney2 (sheet14) ---
$ws = $wb.Worksheets.Item("nomoney2")
$ws.Range("B2").Value = "you being here is enough for now"
$ws.Range("B3").Value = "send whatever feels right, I can't keep this from you"
$ws.Range("B4").Value = "what about just [small amount]? would hate for you to miss this one"
$ws.Range("B5").Value = "don't stress about it"

# --- noppv1 (sheet15) ---
$ws = $wb.Worksheets.Item("noppv1")
$ws.Range("B2").Value = "send whatever you want, literally anything, you need to see this"
$ws.Range("B3").Value = "forget about money, this isn't about that, I need you to see what you're doing to me rn because I've never reacted like this"
$ws.Range("B4").Value = "that's fine, I'm just enjoying this"
$ws.Range("C4").Value = "ACCEPT. Continue sexting 4-5 msgs before Step 2."

# --- noppv2 (sheet16) ---
$ws = $wb.Worksheets.Item("noppv2")
$ws.Range("B2").Value = "send anything, even the smallest amount, you need to see what you did to me"
$ws.Range("B3").Value = "forget the money part, I want to share this with you because what you're making me feel is real"
$ws.Range("B4").Value = "no worries, I don't care about that, this convo is what I care about"
$ws.Range("C4").Value = "ACCEPT. Continue sexting 4-5 msgs before Step 2."

# --- card1 (sheet17) ---
$ws = $wb.Worksheets.Item("card1")
$ws.Range("B2").Value = "figure it out soon because this mood isn't gonna last and I want you to have it"
$ws.Range("B3").Value = "maybe try a different card? you really don't want to miss this"
$ws.Range("B4").Value = "ugh that's annoying, try again it usually works the second time"

# --- card2 (sheet18) ---
$ws = $wb.Worksheets.Item("card2")
$ws.Range("B2").Value = "fix it quick, I don't keep stuff like this around forever"
$ws.Range("B3").Value = "you have another card? because you don't want to miss what I made for you"
$ws.Range("B4").Value = "that happens sometimes, try one more time"

# --- nosex1 (sheet19) ---
$ws = $wb.Worksheets.Item("nosex1")
$ws.Range("B2").Value = "alright I'll behave, for now, no promises"
$ws.Range("B3").Value = "I can't help it though, there's something about you that keeps messing with my head"
$ws.Range("B4").Value = "so tell me more about you, what do you do when you're not driving guys on the internet crazy?"
$ws.Range("B5").Value = "haha okay I got carried away, you're too easy to talk to"
$ws.Range("C2").Value = "ACCEPT. Door open."

# --- nosex2 (sheet20) ---
$ws = $wb.Worksheets.Item("nosex2")
$ws.Range("B2").Value = "fine I'll stop, don't blame me when it happens again though"
$ws.Range("B3").Value = "I'm trying so hard to behave but you make it impossible honestly"
$ws.Range("B4").Value = "okay new topic, what's the craziest thing you've ever done?"
$ws.Range("B5").Value = "my bad I got ahead of myself, your fault for being this fun"

# --- offtopic1 (sheet21) ---
$ws = $wb.Worksheets.Item("offtopic1")
$ws.Range("B2").Value = "okay wait I remember, so like I was saying..."
$ws.Range("B3").Value = "but wait you totally distracted me, I was about to say something and now you made me lose my thought"
$ws.Range("B4").Value = "haha okay that's actually funny"
$ws.Range("C2").Value = "RETAKE. Resume main script."
$ws.Range("C4").Value = "ACKNOWLEDGE. Adapt. → Step 2."

# --- offtopic2 (sheet22) ---
$ws = $wb.Worksheets.Item("offtopic2")
$ws.Range("B2").Value = "OKAY focus, where was I... oh yeah"
$ws.Range("B3").Value = "wait no stop, you're distracting me from what I was gonna tell you"
$ws.Range("B4").Value = "okay that's random but I'm here for it"
$ws.Range("C4").Value = "ACKNOWLEDGE. Adapt. → Step 2."

# --- real1 (sheet23) ---
$ws = $wb.Worksheets.Item("real1")
$ws.Range("B2").Value = "I get why you'd think that, there's a lot of fake stuff on here, but what we've been talking about felt real to me, didn't it feel real to you?"
$ws.Range("B3").Value = "ask me anything, literally anything about me or my life, I'll prove it, go ahead"
$ws.Range("B4").Value = "lol a robot? do robots make you feel like this? didn't think so"

# --- real2 (sheet24) ---
$ws = $wb.Worksheets.Item("real2")
$ws.Range("B2").Value = "I know there's a lot of bots on here but what we've been talking about felt different, you felt it too right?"
$ws.Range("B3").Value = "test me then, ask me something only a real person would know"
$ws.Range("B4").Value = "wait you think I'm fake?? that's the funniest thing I've heard all day"

# --- voice1 (sheet25) ---
$ws = $wb.Worksheets.Item("voice1")
$ws.Range("B2").Value = "I don't do that on here but what I'm about to show you is better than any call"
$ws.Range("B3").Value = "I have something way better for you, trust me you'll forget you even asked"
$ws.Range("B4").Value = "haha maybe if you earn it, not yet though"
$ws.Range("C2").Value = "FIRM. No video calls."
$ws.Range("C4").Value = "DODGE. Still → Step 2."

# --- voice2 (sheet26) ---
$ws = $wb.Worksheets.Item("voice2")
$ws.Range("B2").Value = "that's not something I do here but trust me what I have is way better"
$ws.Range("B3").Value = "how about instead of a call I show you something that'll blow your mind?"
$ws.Range("B4").Value = "you gotta earn that first"

# --- customyes1 (sheet27) ---
$ws = $wb.Worksheets.Item("customyes1")
$ws.Range("B2").Value = "trust me you won't regret it, this one is special"
$ws.Range("B3").Value = "I have exactly what you need and you're gonna lose your mind, [price]"
$ws.Range("B4").Value = "mmm you want that? I might have exactly what you're thinking of"
$ws.Range("C3").Value = "PRICE. Set based on content."

# --- customyes2 (sheet28) ---
$ws = $wb.Worksheets.Item("customyes2")
$ws.Range("B2").Value = "you're not gonna be able to stop watching"
$ws.Range("B3").Value = "I have it, [price] and it's worth every penny"
$ws.Range("B4").Value = "ohhh good taste, I think I know exactly what you need"

# --- customno1 (sheet29) ---
$ws = $wb.Worksheets.Item("customno1")
$ws.Range("B2").Value = "trust me, I know what you need better than you do"
$ws.Range("B3").Value = "what I have might be even crazier and nobody else has seen it"
$ws.Range("B4").Value = "I don't have exactly that but I have something that'll make you forget you asked"
$ws.Range("C3").Value = "ALTERNATIVE + FOMO. → Step 3."

# --- customno2 (sheet30) ---
$ws = $wb.Worksheets.Item("customno2")
$ws.Range("B2").Value = "just trust me on this, you'll thank me after"
$ws.Range("B3").Value = "nobody has seen what I'm about to show you and it's better than what you asked for"
$ws.Range("B4").Value = "not exactly that but what I DO have is gonna hit even harder"

# --- done1 (sheet31) ---
$ws = $wb.Worksheets.Item("done1")
$ws.Range("B2").Value = "next time you have to wait for me, I have something insane planned for round 2"
$ws.Range("B3").Value = "but I haven't finished yet, you're really gonna leave me like this?"
$ws.Range("B4").Value = "fuck that's hot, you came because of me??"

# --- done2 (sheet32) ---
$ws = $wb.Worksheets.Item("done2")
$ws.Range("B2").Value = "next time you hold it because what I have planned is way crazier"
$ws.Range("B3").Value = "wait I'm not done yet, you're just gonna leave me hanging?"
$ws.Range("B4").Value = "already?? damn that's hot"

# --- cumcontrol1 (sheet33) ---
$ws = $wb.Worksheets.Item("cumcontrol1")
$ws.Range("B2").Value = "don't you dare finish before you see this"
$ws.Range("B3").Value = "hold it, wait until you see what I'm about to send, trust me it's worth it"
$ws.Range("B4").Value = "wait for me, I want us to finish together, open this first"
$ws.Range("B5").Value = "I'm close too, cum with me, but you need to see this first"
$ws.Range("B6").Value = "hold it, not yet, I need you to last longer for me"
$ws.Range("B7").Value = "don't cum yet, I'm not done with you"
$ws.Range("C6").Value = "EDGE variant."

# --- dickpic (sheet34) ---
$ws = $wb.Worksheets.Item("dickpic")
$ws.Range("B2").Value = "okay you just made me do something, give me a sec"
$ws.Range("B3").Value = "you can't send me that and expect me to do nothing about it, hold on"
$ws.Range("B4").Value = "woah wasn't expecting that but... damn"
$ws.Range("B5").Value = "omg you don't waste time huh, that's actually really hot ngl"
$ws.Range("B6").Value = "damn that is... I need to show you something rn"
$ws.Range("B7").Value = "fuck okay that's... wow, you have no idea what that just did to me"
$ws.Range("C3").Value = "LEVERAGE. WAIT 1-2 min then send PPV."

# --- boosters (sheet35) ---
$ws = $wb.Worksheets.Item("boosters")
$ws.Range("B2").Value = "I can't focus on anything else rn"
$ws.Range("B3").Value = "more"
$ws.Range("B4").Value = "my hands are shaking"
$ws.Range("B5").Value = "I can't think straight rn"
$ws.Range("B6").Value = "you have no idea what you're doing to me"
$ws.Range("B8").Value = "I'm so hard rn because of you"
$ws.Range("C2").Value = "BOOSTER."
$ws.Range("C9").Value = "MID-SEXTING BOOSTER."

# --- rename sheet33: cumcontrol1 -> cumcontrol ---
$wb.Worksheets.Item("cumcontrol1").Name = "cumcontrol"

